$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header strings & other labels, written in an order that reproduces
# ---- the target sharedStrings sequence: Source, UI, API, CelsiusTemp,
# ---- Mean, x-u, (x-u)^2, Sum, Variance, VarianceLogic
$ws.Range("A1").Value = "Source"
$ws.Range("A2").Value = "UI"
$ws.Range("A3").Value = "API"
$ws.Range("B1").Value = "CelsiusTemp"
$ws.Range("A5").Value = "Mean"
$ws.Range("C1").Value = "x-u"
$ws.Range("D1").Value = "(x-u)^2"
$ws.Range("C5").Value = "Sum"
$ws.Range("C6").Value = "Variance"
$ws.Range("C7").Value = "VarianceLogic"

# ---- Raw numeric inputs
$ws.Range("B2").Value = 30
$ws.Range("B3").Value = 35.89

# ---- Formulas
$ws.Range("C2").Formula = "=B2-B5"
$ws.Range("D2").Formula = "=C2^2"
$ws.Range("C3").Formula = "=B3-B5"
$ws.Range("D3").Formula = "=C3^2"
$ws.Range("B5").Formula = "=AVERAGE(B2:B3)"
$ws.Range("D5").Formula = "=D2+D3"
$ws.Range("D6").Formula = "=D5/2"
$ws.Range("D7").Formula = "=IF(D6<3, TRUE, FALSE)"

# ---- Formatting ----
# Body cells (rows 2-8) get the plain thin-bordered style already used by
# A2:C8; propagate it onto the new column D (and any gaps) via format copy
# from an existing bordered cell so no redundant style entries are minted.
$ws.Range("A2").Copy()
$ws.Range("D2:D8").PasteSpecial(-4122)
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Range("A6:C8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row: A1 already carries the bold / fill / border style. Re-use it
# (format copy) for B1, C1 and D1, then layer on alignment / bold tweaks.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108
$ws.Range("D1").Font.Bold = $false

# ---- Column widths (best effort; engine quantizes to 1/6 char units) ----
$ws.Columns("A:A").ColumnWidth = 13.166666666666666
$ws.Columns("B:B").ColumnWidth = 11.666666666666666
$ws.Columns("C:C").ColumnWidth = 13.333333333333334
$ws.Columns("D:D").ColumnWidth = 13.333333333333334

# ---- Selection ----
$ws.Range("B12").Select()

Write-Output "done"
